$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Range("J1").Value = "target_image"

# New values for column J, rows 2-22
$targetImages = @{
    2  = "stim3"
    3  = "stim2"
    4  = "stim4"
    5  = "stim1"
    6  = "stim3"
    7  = "stim1"
    8  = "stim2"
    9  = "stim1"
    10 = "stim4"
    11 = "simt2"
    12 = "stim3"
    13 = "stim1"
    14 = "stim2"
    15 = "stim4"
    16 = "stim3"
    17 = "stim2"
    18 = "stim4"
    19 = "stim1"
    20 = "stim4"
    21 = "stim3"
    22 = "stim2"
}

# Rows that get the shaded formatting already used elsewhere in the sheet
# (same cell format as e.g. D2, which already carries that style).
$styledRows = @(2,4,6,8,10,11,12,13,14,15,16,17,18,19,20,21,22)
$styleTemplate = $ws.Cells.Item(2, 4)  # D2 already uses the shaded style

foreach ($r in 2..22) {
    $cell = $ws.Cells.Item($r, 10)
    if ($styledRows -contains $r) {
        $styleTemplate.Copy() | Out-Null
        $cell.PasteSpecial(-4122) | Out-Null
    }
    $cell.Value = $targetImages[$r]
}

$excel.CutCopyMode = 0

# Update selection to mirror the edited workbook state
$ws.Range("I12").Select()
